$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with the new "P2209-00002" record ---
$ws.Range("A2").Value = "P2209-00002"
$ws.Range("C2").Value = "P2209-00002"
$ws.Range("D2").Value = "27-09-2022 00:00:00"
$ws.Range("E2").Value = "900.000 VND"
$ws.Range("F2").Value = "P2209-00002"

# --- Remove the now-obsolete rows 3-5 ---
$ws.Range("A3:A5").EntireRow.Delete($null)

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 17 - (5/7)
$ws.Columns.Item(6).ColumnWidth = 13 - (5/7)

# --- Widen the used range out to column V (18-22) with the standard 9.10 width ---
for ($i = 18; $i -le 22; $i++) {
    $ws.Columns.Item($i).ColumnWidth = 9.1 - (5/7)
}

# --- Bump the default font size for the sheet from 11 to 14 ---
$wb.Styles.Item("Normal").Font.Size = 14

Write-Host "done"
